$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, pushing existing rows 21-33 down to 22-34.
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with the new weekly price record.
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 45001
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 100112044
$ws.Range("G21").Value = "Perejil"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 1300
$ws.Range("L21").Value = 1300
$ws.Range("M21").Value = 1300
$ws.Range("N21").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 1300
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"
